$d = $word.ActiveDocument
$t = $d.Tables(1)

$replacements = @(
    @{Row=1;  Col=1; Old="52÷7="; New="19÷4="},
    @{Row=1;  Col=2; Old="90÷9="; New="81÷7="},
    @{Row=1;  Col=3; Old="39÷3="; New="70÷5="},
    @{Row=1;  Col=4; Old="48÷5="; New="22÷7="},
    @{Row=1;  Col=5; Old="77÷2="; New="54÷6="},

    @{Row=5;  Col=1; Old="39÷3="; New="61÷8="},
    @{Row=5;  Col=2; Old="63÷3="; New="58÷3="},
    @{Row=5;  Col=3; Old="70÷6="; New="57÷2="},
    @{Row=5;  Col=4; Old="36÷6="; New="35÷8="},
    @{Row=5;  Col=5; Old="42÷5="; New="71÷2="},

    @{Row=9;  Col=1; Old="84÷4="; New="28÷6="},
    @{Row=9;  Col=2; Old="65÷6="; New="61÷5="},
    @{Row=9;  Col=3; Old="84÷7="; New="16÷4="},
    @{Row=9;  Col=4; Old="75÷2="; New="31÷8="},
    @{Row=9;  Col=5; Old="83÷5="; New="54÷2="},

    @{Row=13; Col=1; Old="25÷9="; New="60÷2="},
    @{Row=13; Col=2; Old="57÷3="; New="91÷9="},
    @{Row=13; Col=3; Old="12÷3="; New="48÷5="},
    @{Row=13; Col=4; Old="26÷3="; New="54÷6="},
    @{Row=13; Col=5; Old="76÷3="; New="80÷4="},

    @{Row=17; Col=1; Old="81÷2="; New="11÷5="},
    @{Row=17; Col=2; Old="93÷7="; New="16÷7="},
    @{Row=17; Col=3; Old="84÷5="; New="59÷8="},
    @{Row=17; Col=4; Old="77÷6="; New="26÷8="},
    @{Row=17; Col=5; Old="80÷8="; New="88÷3="}
)

foreach ($item in $replacements) {
    $cell = $t.Cell($item.Row, $item.Col)
    $current = $cell.Range.Text
    if (-not $current.StartsWith($item.Old)) {
        Write-Host ("Unexpected text at row " + $item.Row + ", col " + $item.Col + ": [" + $current + "] expected [" + $item.Old + "]")
    }
    $cell.Range.Text = $item.New
}
